$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the paragraph that ends with "LOM3062: Trabalho de Graduação I (Requisito)".
# The three paragraphs that immediately follow it (an empty paragraph, the
# "Ver no Jupiter..." paragraph and the "© 2020 ..." copyright paragraph)
# need to be removed entirely, while the paragraph after that (the blank
# paragraph right before the final page-break paragraph) must stay.
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "LOM3062:*Requisito)*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $startPara = $paras.Item($targetIndex + 1)
    $endPara = $paras.Item($targetIndex + 3)
    $removeRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $removeRange.Delete()
}
